$wb = $excel.ActiveWorkbook

# ---- Create_Contact sheet ----
$wsContact = $wb.Worksheets.Item("Create_Contact")
$wsContact.Activate()

$wsContact.Range("C2").Value = "1Prime"
$wsContact.Range("C5").Value = "1Youtube"
$wsContact.Range("C8").Value = "1whatsapp"

$wsContact.Range("C8").Select()

# ---- Create_Campaign sheet ----
$wsCampaign = $wb.Worksheets.Item("Create_Campaign")
$wsCampaign.Activate()

$wsCampaign.Range("C2").Value = "HiCampaign"
$wsCampaign.Range("D2").Value = "'21"

$wsCampaign.Range("C5").Value = "BB_CampaignWithStatus"
$wsCampaign.Range("E5").Value = "'27"

$wsCampaign.Range("C8").Value = "BB_CampaignWithExpectedClosedate"
$wsCampaign.Range("E8").Value = "'18"

$wsCampaign.Range("C11").Value = "BB_CreateCampaignCompleteTest"
$wsCampaign.Range("E11").Value = "'11"

$wsCampaign.Range("D12").Select()

# ---- restore the active sheet/tab (Create_Contact is tab index 1, activeTab=1) ----
$wsContact.Activate()
